# --------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1. Insert a brand-new "2022-Q3" worksheet right after "总计" (becomes
#    the 2nd tab), with the same per-fund layout as the other quarters.
# 2. Insert the 2022-Q3 totals as the first data row of "总计", shifting
#    every other quarter's row down by one and renumbering the index
#    column.
# --------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. New "2022-Q3" sheet -------------------------------------------------
$anchor    = $wb.Worksheets.Item(2)
$newSheet  = $wb.Worksheets.Add($anchor)
$newSheet.Name = "2022-Q3"

# Reuse the header / index-column formatting from an existing quarterly
# sheet so the new tab looks like the rest of the workbook.
$template = $wb.Worksheets.Item(3)
$template.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$template.Range("A2:A3").Copy($newSheet.Range("A2:A4"))

# Header row (matches the other quarterly fund-holding sheets).
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Fund code / amount columns are kept as text (so "004818", "2.70",
# "45.00", ... keep their original printed form instead of becoming
# numbers that drop leading/trailing zeros).
$newSheet.Range("B2:G4").NumberFormat = "@"

$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "004818"
$newSheet.Cells.Item(2,3).Value = "国寿安保目标策略灵活配置混合A"
$newSheet.Cells.Item(2,4).Value = "2.70"
$newSheet.Cells.Item(2,5).Value = "45.00"
$newSheet.Cells.Item(2,6).Value = "2.09"
$newSheet.Cells.Item(2,7).Value = "0.0564"
$newSheet.Cells.Item(2,8).Value = 8

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "004819"
$newSheet.Cells.Item(3,3).Value = "国寿安保目标策略灵活配置混合C"
$newSheet.Cells.Item(3,4).Value = "1.73"
$newSheet.Cells.Item(3,5).Value = "45.00"
$newSheet.Cells.Item(3,6).Value = "2.09"
$newSheet.Cells.Item(3,7).Value = "0.0362"
$newSheet.Cells.Item(3,8).Value = 8

$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "003238"
$newSheet.Cells.Item(4,3).Value = "新华外延增长主题灵活配置混合"
$newSheet.Cells.Item(4,4).Value = "0.50"
$newSheet.Cells.Item(4,5).Value = "57.43"
$newSheet.Cells.Item(4,6).Value = "3.09"
$newSheet.Cells.Item(4,7).Value = "0.0154"
$newSheet.Cells.Item(4,8).Value = 3

# ---- 2. "总计" summary sheet: push rows down, add 2022-Q3 on top -----------
$summary = $wb.Worksheets.Item(1)

# Copy the formatting of the last existing data row down into the new
# row 7 before writing to it, so the index column keeps its styling.
$summary.Cells.Item(6,1).Copy($summary.Cells.Item(7,1))

$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(7,2).Value = "2020-Q4"
$summary.Cells.Item(7,3).Value = 7
$summary.Cells.Item(7,4).Value = 5.05

$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = "2021-Q1"
$summary.Cells.Item(6,3).Value = 8
$summary.Cells.Item(6,4).Value = 6.67

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q2"
$summary.Cells.Item(5,3).Value = 1
$summary.Cells.Item(5,4).Value = 0.01

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2021-Q4"
$summary.Cells.Item(4,3).Value = 4
$summary.Cells.Item(4,4).Value = 3.48

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q1"
$summary.Cells.Item(3,3).Value = 2
$summary.Cells.Item(3,4).Value = 0.96

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 0.11
